$wb = $excel.ActiveWorkbook

# --- Update the "Conversion del dia" text on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.79 = 6472.27 pesos`n✅ 6472.27 pesos = 1.78 = 935.68 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate values on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 559
$ws2.Range("O10").Value = 3618
$ws2.Range("N12").Value = 3636
$ws2.Range("O12").Value = 525.65
